# SignalMap_KPI_PlotSpec.xlsx - "load aeb relevant params from class and config"
#
# The "params" sheet gains two new AEB-specific timing parameters
# (PRE_TIME_AEB / POST_TIME_AEB) grouped with the other AebEventDetector
# rows, loses the old FCW-specific timing parameters (PRE_TIME_FCW /
# POST_TIME_FCW, class FcwEventDetector) which are superseded by new
# generic defaults (PRE_TIME_DEFAULT / POST_TIME_DEFAULT, class
# BaseEventSegmenter) appended at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")

# --- Insert two new rows right above the existing AebEventDetector rows
# (old row 4, "START_DECEL_DELTA") so PRE_TIME_AEB / POST_TIME_AEB lead
# that class's block. ---
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

$ws.Range("A4").Value = "PRE_TIME_AEB"
$ws.Range("B4").Value = 6
$ws.Range("C4").Value = "float"
$ws.Range("D4").Value = "s"
$ws.Range("E4").Value = "time before event (duration)"
$ws.Range("F4").Value = "AebEventDetector"

$ws.Range("A5").Value = "POST_TIME_AEB"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "float"
$ws.Range("D5").Value = "s"
$ws.Range("E5").Value = "time after event (duration)"
$ws.Range("F5").Value = "AebEventDetector"

# --- Remove the old FCW-specific PRE_TIME_FCW / POST_TIME_FCW rows (now
# shifted down to rows 13/14 by the insert above). The FcwEventDetector
# class they referenced is being replaced by BaseEventSegmenter. ---
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(13).Delete()

# --- Append the new generic default timing parameters at the bottom of
# the table (rows 18/19). ---
$ws.Range("A18").Value = "PRE_TIME_DEFAULT"
$ws.Range("B18").Value = 6
$ws.Range("C18").Value = "float"
$ws.Range("D18").Value = "s"
$ws.Range("E18").Value = "default time before event (duration)"
$ws.Range("F18").Value = "BaseEventSegmenter"

$ws.Range("A19").Value = "POST_TIME_DEFAULT"
$ws.Range("B19").Value = 3
$ws.Range("C19").Value = "float"
$ws.Range("D19").Value = "s"
$ws.Range("E19").Value = "default time after event (duration)"
$ws.Range("F19").Value = "BaseEventSegmenter"

# --- Restore the active selection on the "params" sheet (was sitting on
# a stale E25 selection from when the sheet was shorter). ---
$ws.Activate()
$ws.Range("E7").Select()
